$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    $r = $ws.Range($cellAddr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" '37.188.99'
$ws.Range("E2").Value = '  +1.45%  '
Set-TextValue "D3" '2.024.73'
$ws.Range("E3").Value = '  +3.12%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("E6").Value = '  +1.80%  '
Set-TextValue "D7" '60.48'
$ws.Range("E7").Value = '  -1.97%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +4.70%  '
$ws.Range("E10").Value = '  +2.09%  '
$ws.Range("E11").Value = '  +1.88%  '
Set-TextValue "D12" '15.22'
$ws.Range("E12").Value = '  +6.24%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D13" '0.860'
$ws.Range("E13").Value = '  +3.00%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D14" '22.50'
$ws.Range("E14").Value = '  +1.53%  '
Set-TextValue "D15" '2.322.75'
$ws.Range("E15").Value = '  +3.21%  '
Set-TextValue "D16" '5.53'
$ws.Range("E16").Value = '  +4.11%  '
Set-TextValue "D17" '2.024.81'
$ws.Range("E17").Value = '  +2.85%  '
Set-TextValue "D18" '37.149.22'
$ws.Range("E18").Value = '  +1.38%  '
Set-TextValue "D19" '70.69'
$ws.Range("E19").Value = '  +1.18%  '
Set-TextValue "D20" '0.0₃0868'
$ws.Range("E20").Value = '  +1.35%  '
Set-TextValue "D21" '5.26'
$ws.Range("E21").Value = '  +3.22%  '
Set-TextValue "D22" '231.36'
$ws.Range("E22").Value = '  +0.42%  '
Set-TextValue "D24" '2.51'
$ws.Range("E24").Value = '  +2.07%  '
$ws.Range("E25").Value = '  +0.70%  '
Set-TextValue "D26" '9.47'
$ws.Range("E26").Value = '  +2.43%  '
Set-TextValue "D27" '164.01'
$ws.Range("E27").Value = '  +1.97%  '
Set-TextValue "D28" '0.139'
$ws.Range("E28").Value = '  -3.22%  '
Set-TextValue "D29" '19.85'
$ws.Range("E29").Value = '  +2.08%  '
$ws.Range("E30").Value = '  +7.62%  '
$ws.Range("E31").Value = '  +2.01%  '
$ws.Range("E32").Value = '  +1.12%  '
Set-TextValue "D33" '0.0667'
Set-TextValue "D34" '4.57'
$ws.Range("E34").Value = '  +2.33%  '
$ws.Range("E35").Value = '  +9.61%  '
$ws.Range("E36").Value = '  -3.78%  '
$ws.Range("E37").Value = '  -0.12%  '
Set-TextValue "D39" '5.45'
$ws.Range("E39").Value = '  -1.57%  '
Set-TextValue "D40" '0.0982'
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  +1.22%  '
Set-TextValue "D42" '1.20'
$ws.Range("E42").Value = '  +1.34%  '
$ws.Range("E43").Value = '  +1.56%  '
Set-TextValue "D44" '16.88'
$ws.Range("E44").Value = '  +4.75%  '
Set-TextValue "D45" '92.27'
$ws.Range("E45").Value = '  +3.76%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D46" '1.390.35'
$ws.Range("E46").Value = '  +1.38%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D47" '1.07'
$ws.Range("E47").Value = '  +3.12%  '
Set-TextValue "D48" '7.51'
$ws.Range("E48").Value = '  +4.83%  '
Set-TextValue "D49" '2.17'
$ws.Range("E49").Value = '  +17.74%  '
$ws.Range("E50").Value = '  +0.42%  '
Set-TextValue "D51" '46.86'
$ws.Range("E51").Value = '  +3.07%  '
